# Applies the "confronto-pratiche-Progetto-Angelo" update:
#  - recompute the per-locator-type failure table (rows 11-16): the
#    "Totale" column (D) becomes a plain value instead of a formula,
#    while "Fallimenti per Fragilita'" (E) and "Fallimenti per
#    Obsolescenza" (F) get new counts; the "Tasso di fallimento" (G)
#    formula is left in place and recalculates automatically.
#  - row 12's D cell loses its (now unused) underlined-font formatting.
#  - the totals row (21) formulas recalculate automatically from the
#    updated E/F columns.
#  - the sheet's saved selection moves from E21 to B21 and the frozen
#    top-left cell (B1) is cleared.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 11 (hook) ---
$ws.Range("D11").Value = 30
$ws.Range("E11").Value = 5
$ws.Range("F11").Value = 3

# --- Row 12 (absolute) ---
$ws.Range("D12").Value = 20
$ws.Range("E12").Value = 15
$ws.Range("F12").Value = 3
# Drop the underline formatting this cell used to carry (its style
# becomes equivalent to the plain centered style used elsewhere).
$ws.Range("D12").Font.Underline = -4142

# --- Row 13 (relative) ---
$ws.Range("D13").Value = 29
$ws.Range("E13").Value = 6
$ws.Range("F13").Value = 3

# --- Row 14 (robula) ---
$ws.Range("D14").Value = 33
$ws.Range("E14").Value = 2
$ws.Range("F14").Value = 3

# --- Row 15 (selenium) ---
$ws.Range("D15").Value = 29
$ws.Range("E15").Value = 6
$ws.Range("F15").Value = 3

# --- Row 16 (katalon) ---
$ws.Range("D16").Value = 33
$ws.Range("E16").Value = 2
$ws.Range("F16").Value = 3

# Row 21 totals (B21=SUM(E11:E16), C21=SUM(F11:F16), D21=SUM(B21,C21))
# are formulas already present on the sheet, so they recompute on
# their own from the edits above.

# Update the saved view: select B21 (instead of E21) and let Excel
# drop the explicit top-left-cell freeze that used to pin the view at B1.
$ws.Activate()
$ws.Range("B21").Select()
